$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing cells: B39 and B40
$ws.Range("B39").Value = 464
$ws.Range("B40").Value = 522

# Add new rows 41 and 42
$ws.Range("A41").Value = 40
$ws.Range("B41").Value = 424

$ws.Range("A42").Value = 41
$ws.Range("B42").Value = 23
